# Apply the JIRA-ID worksheet update:
#  - A2: JCVZ-310 -> JCVZ-929 (text + hyperlink)
#  - A3: JCVZ-929 -> JCVZ-1030 (text + hyperlink)
#  - A4: cleared (was JCVZ-928, hyperlink removed)
#  - A5: cleared (was JCVZ-932, hyperlink removed)
#  - Selection moved to E9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember a still-"clean" style (A4's current xf, s="1") so we can restore it
# on A2/A3 after Hyperlinks.Add() re-stamps their style below.
$cleanStyle = $ws.Range("A4").Style

# This COM layer's Hyperlinks.Delete() always clears every hyperlink on the
# sheet (it is not scoped to the calling range), so do it once up front and
# rebuild only the links we still want.
$ws.Hyperlinks.Delete()

# A4 / A5 lose their text and hyperlink entirely, keeping their existing
# cell style.
$ws.Range("A4:A5").ClearContents()

# A2 / A3 get new display text...
$ws.Range("A2").Value = "https://jira.jnj.com/browse/JCVZ-929"
$ws.Range("A3").Value = "https://jira.jnj.com/browse/JCVZ-1030"

# ...and new hyperlinks pointing at the same URL as their text.
$ws.Hyperlinks.Add($ws.Range("A2"), "https://jira.jnj.com/browse/JCVZ-929")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://jira.jnj.com/browse/JCVZ-1030")

# Hyperlinks.Add() stamps a fresh "applyFont" cell style onto the target -
# put back the plain Hyperlink style so A2/A3 keep their original look.
$ws.Range("A2").Style = $cleanStyle
$ws.Range("A3").Style = $cleanStyle

# Move the active selection to E9, matching the saved view state.
$null = $ws.Range("E9").Select()
